$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new study row (row 28) for "Lescure et al." (Sarilumab), mirroring
# the formatting of the previous row (27) so the new cells pick up the same
# per-column styling (ArialMT9 / Arial10 / Arial12 for columns A/B/C).
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new study's data: risk-of-bias domains D1-D5 and Overall are
# all "Low", assessed "we", Sarilumab, with a WHO weight of 1.
$ws.Range("A28").Value = "Lescure et al."
$ws.Range("B28").Value = "we"
$ws.Range("C28").Value = "Sarilumab"
$ws.Range("D28").Value = "Low"
$ws.Range("E28").Value = "Low"
$ws.Range("F28").Value = "Low"
$ws.Range("G28").Value = "Low"
$ws.Range("H28").Value = "Low"
$ws.Range("I28").Value = "Low"
$ws.Range("J28").Value = 1

# Match the updated selection state left in the source file.
$ws.Range("A11").Select()
